# --------------------------------------------------------------------------
# "error solve ifrs list" - fixes a unit/scale bug in the IFRS company_list
# sheet:
#   - rows 2-6 (FY2014-FY2018 actuals) had every financial metric (D:AJ) off by
#     a wrong scale/shift; they are rewritten here with the corrected figures.
#   - rows 7-9 (FY2019(E)-FY2021(E) estimates) are dropped entirely: everything
#     except the A/B/C label columns is cleared, matching upstream removing
#     those (unreliable) estimate cells.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ")
$vals = @(691, 16, 16, -34, -40, -40, 0, 1955, 1229, 726, 731, -5, 181, -70, -17, 95, 3, -73, 880, 2.29, -5.76, -5.21, -2.01, 169.15, 354.21, -109, -10.44, 3470, 0.33, 25, 2.19, -9.460000000000001, 36212160)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $vals[$i]
}

# Row 3
$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ")
$vals = @(759, 28, 28, -3, -5, -4, -1, 2017, 1307, 710, 716, -6, 181, -81, 2, 74, 0, -82, 975, 3.67, -0.6899999999999999, -0.6, -0.26, 184.16, 352.8, -12, -166.17, 3398, 0.59, 25, 1.25, -36.17, 36212160)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $vals[$i]
}

# Row 4
$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ")
$vals = @(837, 46, 46, 17, 13, 14, 0, 2045, 1312, 734, 741, -7, 181, 69, 37, -98, 3, 65, 894, 5.45, 1.59, 1.87, 0.65, 178.71, 367.26, 38, 48.85, 3320, 0.55, 25, 1.36, 14.67, 36212160)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $vals[$i]
}

# Row 5
$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ")
$vals = @(920, 72, 72, 33, 14, 14, 0, 2034, 1160, 874, 876, -2, 181, 67, -11, -47, 2, 66, 716, 7.84, 1.47, 1.68, 0.66, 132.78, 409.53, 37, 57.73, 3206, 0.67, 25, 1.16, 23.62, 36212160)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $vals[$i]
}

# Row 6
$cols = @("D", "E", "F", "G", "H", "I", "K", "L", "M", "N", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ")
$vals = @(833, 10, 10, -28, -36, -35, 2051, 1223, 828, 831, 181, 10, -43, 52, 2, 8, 792, 1.17, -4.3, -4.11, -1.75, 147.78, 385.27, -97, -19.73, 3055, 0.63, 25, 1.31, -8.98, 36212160)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $vals[$i]
}

# Rows 7-9: the FY2019(E)-FY2021(E) estimate columns are removed outright,
# leaving only the A (index), B (period) and C (year) label cells behind.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
